# Add two new columns (I: "I0", J: "IF") to Sheet1, mirroring the header
# style already used by the other header cells (B1:H1), then fill in the
# per-row values for rows 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting from the existing "IP" header (H1) onto the two new
# header cells so they share the same bold/border/alignment style, then
# set their text.
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 10).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data ----------------------------------------------------------------
$I0 = @(7, 8, 8, 5, 6, 6, 1, 8, 7, 7, 6, 6, 8, 7, 8, 6, 7, 7, 8, 10, 6, 8, 7, 6, 2, 6, 3, 9, 6, 7, 4, 6, 3)
$IF = @(7, 8, 8, 5, 7, 6, 3, 8, 7, 7, 6, 6, 8, 7, 8, 7, 7, 7, 8, 10, 7, 8, 7, 7, 2, 6, 4, 9, 6, 7, 5, 6, 3)

for ($i = 0; $i -lt $I0.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
